# netCrypto.xlsx update — user corrected the USD Amount figure in T2 and
# then moved the selection down to T3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the USD Amount value in T2 (105329 -> 103842)
$ws.Range("T2").Value = 103842

# Leave the active selection on T3, matching where the user clicked next
$ws.Range("T3").Select()
